# Add team record (Wins / Losses / Ties) columns to the KCR_1991 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1=Wins, AE1=Losses, AF1=Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, bordered, centered) from AC1
# onto the new header cells, the same way the other header columns look.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

# Data rows 2..47: every player row gets the team's 1991 Royals record (82-80-0)
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 82   # column AD -> Wins
    $ws.Cells.Item($r, 31).Value = 80   # column AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF -> Ties
}
